# Update the division problems in the practice table to the new set of
# three-digit-number / one-digit-number problems.
#
# Each entry replaces the old "XXX÷Y=" expression with the new one. Since
# every "old" string is a unique, literal substring in the document, a
# simple MatchCase Find/Replace (ReplaceAll) for each pair is sufficient
# and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @("774÷6=", "920÷6="),
    @("482÷2=", "357÷4="),
    @("276÷7=", "485÷3="),
    @("287÷8=", "422÷2="),
    @("391÷7=", "882÷6="),
    @("579÷2=", "758÷9="),
    @("411÷5=", "128÷9="),
    @("724÷3=", "766÷3="),
    @("530÷5=", "232÷2="),
    @("424÷4=", "403÷8="),
    @("325÷7=", "886÷8="),
    @("509÷2=", "333÷5="),
    @("962÷3=", "363÷5="),
    @("126÷8=", "678÷3="),
    @("833÷9=", "329÷6="),
    @("938÷9=", "698÷5="),
    @("275÷7=", "115÷4="),
    @("618÷2=", "841÷6="),
    @("393÷3=", "648÷9="),
    @("523÷4=", "402÷2="),
    @("937÷6=", "190÷9="),
    @("538÷4=", "536÷5="),
    @("771÷8=", "866÷5="),
    @("823÷5=", "200÷6="),
    @("686÷5=", "684÷9=")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]

    $d.Content.Find.Execute(
        $oldText,  # FindText
        $true,     # MatchCase
        $false,    # MatchWholeWord
        $false,    # MatchWildcards
        $false,    # MatchSoundsLike
        $false,    # MatchAllWordForms
        $true,     # Forward
        1,         # Wrap (wdFindContinue)
        $false,    # Format
        $newText,  # ReplaceWith
        2          # Replace (wdReplaceAll)
    ) | Out-Null
}
